$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "2003年" row (row 2) and the "2008年" row (row 3).
# Deleting entire rows shifts the remaining data (2013年, 2018年) up.
$ws.Rows("2:3").Delete()
